# Refresh the cryptos price table (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain-text values (e.g. "1.234.56"); a leading
# apostrophe forces text entry so Excel doesn't reinterpret values such as
# "1.002" or "47.31" as numbers.

$ws.Range("D2").Value = "'24.802.63"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "'1.664.83"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'330.73"
$ws.Range("E5").Value = "  +8.29%  "

$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("D8").Value = "'47.31"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'0.3244"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("D11").Value = "'0.07063"
$ws.Range("E11").Value = "  +2.41%  "

$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "'6.071"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").Value = "'19.62"
$ws.Range("E14").Value = "  +2.85%  "

$ws.Range("D15").Value = "'1.663.26"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").Value = "'6.608"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").Value = "'0.00001051"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").Value = "'0.06621"
$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("E20").Value = "  +2.61%  "

$ws.Range("D21").Value = "'5.934"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").Value = "'15.83"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("E23").Value = "  +3.30%  "

$ws.Range("D24").Value = "'24.785.50"
$ws.Range("E24").Value = "  +1.62%  "

$ws.Range("D25").Value = "'2.452"
$ws.Range("E25").Value = "  +2.15%  "

$ws.Range("D26").Value = "'2.412"
$ws.Range("E26").Value = "  +3.77%  "

$ws.Range("D27").Value = "'148.70"
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("D28").Value = "'18.67"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("D29").Value = "'1.846.37"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("D30").Value = "'125.61"
$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("D31").Value = "'1.175"
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").Value = "'5.706"
$ws.Range("E33").Value = "  +2.22%  "

$ws.Range("D34").Value = "'0.08489"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("D35").Value = "'1.638"
$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("D37").Value = "'5.160"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").Value = "'0.02261"
$ws.Range("E38").Value = "  +2.45%  "

$ws.Range("D39").Value = "'0.06046"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.225"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2082"
$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("D42").Value = "'8.211"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("D44").Value = "'0.5923"
$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("D45").Value = "'13.40"
$ws.Range("E45").Value = "  +6.48%  "

$ws.Range("D46").Value = "'3.852"
$ws.Range("E46").Value = "  +3.73%  "

$ws.Range("D47").Value = "'0.5663"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("D48").Value = "'125.51"
$ws.Range("E48").Value = "  +3.45%  "

$ws.Range("D49").Value = "'1.949"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").Value = "'0.06966"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").Value = "'1.186"
$ws.Range("E51").Value = "  +3.95%  "
